# Add the new "ODI Bowling Extra" worksheet after the last existing sheet
# ("ODI Batting Extra"), matching sheetId=5 / tab order position.
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "ODI Bowling Extra"

# Reuse the header formatting (bold, bordered, centered) already used by the
# "ODI Batting Extra" sheet's header row so the new header cells share the
# same style instead of creating a brand-new one.
$srcHeader = $wb.Worksheets.Item("ODI Batting Extra").Range("A1:C1")
$srcHeader.Copy($ws.Range("A1:C1"))

$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "MAIDEN_OVERS"
$ws.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

# Data rows. All values are stored as literal text (matching the source
# data feed), so force the cells to Text format before typing anything that
# looks like a number or a percentage.
$rows = @(
  @{ Code = "4294"; Maiden = $null; Pct = $null },
  @{ Code = "4300"; Maiden = "0";   Pct = "40.00%" },
  @{ Code = "4319"; Maiden = "0";   Pct = "20.00%" },
  @{ Code = "4334"; Maiden = "0";   Pct = "10.00%" },
  @{ Code = "4337"; Maiden = $null; Pct = $null },
  @{ Code = "4340"; Maiden = "0";   Pct = "40.00%" },
  @{ Code = "4349"; Maiden = "0";   Pct = "60.00%" },
  @{ Code = "4432"; Maiden = $null; Pct = $null },
  @{ Code = "4433"; Maiden = $null; Pct = $null },
  @{ Code = "4434"; Maiden = "1";   Pct = $null },
  @{ Code = "4458"; Maiden = "1";   Pct = "20.00%" },
  @{ Code = "4459"; Maiden = "1";   Pct = "10.00%" },
  @{ Code = "4460"; Maiden = "0";   Pct = "30.00%" },
  @{ Code = "4472"; Maiden = "0";   Pct = "10.00%" },
  @{ Code = "4473"; Maiden = "0";   Pct = "10.00%" },
  @{ Code = "4476"; Maiden = "0";   Pct = $null },
  @{ Code = "4565"; Maiden = $null; Pct = $null },
  @{ Code = "4567"; Maiden = "2";   Pct = "20.00%" },
  @{ Code = "4586"; Maiden = $null; Pct = $null },
  @{ Code = "4590"; Maiden = $null; Pct = $null }
)

$r = 2
foreach ($row in $rows) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row.Code

    if ($row.Maiden -ne $null) {
        $cellB = $ws.Cells.Item($r, 2)
        $cellB.NumberFormat = "@"
        $cellB.Value = $row.Maiden
    }
    if ($row.Pct -ne $null) {
        $cellC = $ws.Cells.Item($r, 3)
        $cellC.NumberFormat = "@"
        $cellC.Value = $row.Pct
    }
    $r = $r + 1
}
